$d = $word.ActiveDocument

$d.Content.Find.Execute("Return 47", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Return 23", 2)

$d.Content.Find.Execute("Returns 37", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Returns 12", 2)

$d.Content.Find.Execute("Returns 24", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Returns 0", 2)

$d.Content.Find.Execute("Returns 10", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Returns 9", 2)
